$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: update the order status (column E) to a new "Delivery" value and
# set the driver id (column I) to 8.
$ws.Range("E16").Value = "Delivery"
$ws.Range("I16").Value = 8
